$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.052.24'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.467.85'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.71'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.07'
$ws.Range("E6").Value = '  -2.51%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.467.95'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.96'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.52'
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.831.82'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("E17").Value = '  -3.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.467.26'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.17'
$ws.Range("E19").Value = '  -4.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.56'
$ws.Range("E20").Value = '  -4.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '355.30'
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.13'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("E25").Value = '  -6.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -8.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.594.89'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("E30").Value = '  -5.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '509.45'
$ws.Range("E31").Value = '  -3.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.74'
$ws.Range("E32").Value = '  -7.13%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -5.06%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -4.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("E36").Value = '  -6.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.92'
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.62'
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("E40").Value = '  -4.80%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").Value = '  -5.74%  '
$ws.Range("E43").Value = '  -6.15%  '
$ws.Range("E44").Value = '  -5.81%  '
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.73'
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.80'
$ws.Range("E47").Value = '  -2.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.46'
$ws.Range("E48").Value = '  -5.29%  '
$ws.Range("E49").Value = '  -5.38%  '
$ws.Range("E50").Value = '  -5.14%  '
$ws.Range("E51").Value = '  -8.47%  '
